# Regenerate merged AHB files:
#  - rename the "_old"/"_new" header-name suffixes to "_FV2210"/"_FV2304"
#  - wrap the data range in a named table (Table1)
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cells ------------------------------------
$columns = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $columns.Count; $i++) {
    # columns A..J -> "<name>_FV2210"
    $ws.Cells.Item(1, $i + 1).Value = $columns[$i] + "_FV2210"
    # columns L..U -> "<name>_FV2304" (column K holds "diff", unchanged)
    $ws.Cells.Item(1, $i + 12).Value = $columns[$i] + "_FV2304"
}

# --- 2. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into a table ------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), [System.Type]::Missing, 1, [System.Type]::Missing)
$tbl.Name = "Table1"
